$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.885.49"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.292.74"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.80%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "108.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +11.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "271.41"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.51%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.619"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.45%  "
$ws.Range("E8").Value = "  +0.23%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.613"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.95%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "46.87"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.57%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0935"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.79%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.33"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.107"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.61"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.637.11"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.37%  "
$ws.Range("E16").Value = "  -1.56%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.288.05"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.38%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.819.58"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.22%  "
$ws.Range("E19").Value = "  +0.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.87%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.49"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +9.59%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "233.60"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.46%  "
$ws.Range("E24").Value = "  +15.77%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.31"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.43%  "
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.34"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "40.70"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.75%  "
$ws.Range("E29").Value = "  -0.94%  "
$ws.Range("E30").Value = "  +1.21%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "177.82"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.51%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.86"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.72%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0908"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.64%  "
$ws.Range("E34").Value = "  +1.33%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.89"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +9.84%  "
$ws.Range("E36").Value = "  -0.17%  "
$ws.Range("E37").Value = "  +3.44%  "
$ws.Range("E38").Value = "  -1.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.63"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.59%  "
$ws.Range("E40").Value = "  -3.05%  "
$ws.Range("E41").Value = "  -2.82%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.37"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.79%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "66.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.82%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.27"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.72%  "
$ws.Range("E45").Value = "  +2.49%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.78"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.37%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.101"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.99%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.23"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.68%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "99.58"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.87%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.54"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +11.07%  "
$ws.Range("E51").Value = "  +5.24%  "
